$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as text (e.g. "137.14"); setting NumberFormat to
# "@" (Text) before assigning keeps numeric-looking strings as text instead of
# Excel auto-converting them to Number values.
$dCells = @(
"D2", "D3", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D20", "D22", "D23", "D24", "D25", "D28", "D30", "D33", "D38", "D40", "D42", "D44", "D47"
)
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.980.42"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.418.66"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "137.19"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").Value = "0.106"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "0.148"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "25.34"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").Value = "2.849.08"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "59.917.57"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "2.409.63"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "328.04"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "66.02"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "0.178"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("D25").Value = "8.66"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").Value = "0.0₃0775"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "169.22"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "18.60"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "326.58"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").Value = "0.407"
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("D42").Value = "140.21"
$ws.Range("E42").Value = "  -3.12%  "
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "19.68"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "0.404"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("E51").Value = "  -1.09%  "
